# Added genre and comments controllers
# This script edits the API description sheet:
#  - Removes the "view single genre" row from the GENRE block and renames
#    the block from GENRES to GENRE, updating its endpoints to
#    /genre/add (POST) and /genre/all (GET list).
#  - Removes the "view single subscription" and "edit subscription" rows
#    from the FOLLOWERS block.
#  - Restores the selection to D53 (matching the author's final cursor spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work bottom-to-top so row numbers used below stay valid as rows are removed.

# FOLLOWERS block (rows 52-57 before edit): drop the "edit subscription" (row 55)
# and "view single subscription" (row 54) rows, keeping add/delete/list.
$ws.Rows("55").Delete()
$ws.Rows("54").Delete()

# GENRES block (rows 38-43 before edit): drop the "view single genre" row (row 40).
$ws.Rows("40").Delete()

# Rename the block header GENRES -> GENRE
$ws.Range("A38").Value = "GENRE"

# Update the add-genre endpoint
$ws.Range("C39").Value = "/genre/add"

# Update the list-genres endpoint (row 42 after the row-40 deletion above)
$ws.Range("C42").Value = "/genre/all"

# Restore the selection / active cell shown in the saved workbook.
$ws.Range("D53").Select()
